# GitHub Actions style refresh of the cryptos price table:
# update the "Price" (D) and "Volume(1h)" (E) columns in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.915.14"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.645.83"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "216.28"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "0.5063"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "0.2586"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'0.06440"
$ws.Range("D10").Value = "20.42"
$ws.Range("E10").Value = "  +4.59%  "
$ws.Range("D11").Value = "0.07819"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "4.275"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "1.645.55"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "1.871.88"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "0.5632"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "0.0₅7705"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "'63.40"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "25.922.13"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "193.03"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").Value = "4.377"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").Value = "9.949"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "6.128"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "1.804"
$ws.Range("E25").Value = "  -5.85%  "
$ws.Range("D26").Value = "141.76"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").Value = "0.1239"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").Value = "6.799"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "1.248"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("D31").Value = "0.04946"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").Value = "'3.240"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").Value = "1.577"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("D35").Value = "2.394"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").Value = "0.9067"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").Value = "0.5566"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "1.133.72"
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("D39").Value = "2.548"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "0.01565"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "5.519"
$ws.Range("D43").Value = "0.8047"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").Value = "99.09"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("D45").Value = "1.782.19"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").Value = "  -6.89%  "
$ws.Range("D47").Value = "55.81"
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("D48").Value = "0.4298"
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("D49").Value = "7.741"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").Value = "0.05044"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").Value = "0.9996"
$ws.Range("E51").Value = "  -0.56%  "
